# Quarterly indexing esoteric bug-fix operation
#
# For each data row (2..16) the per-quarter error series needs to be
# re-indexed: a newly computed value is inserted at column B and every
# existing value (B..J) shifts one column to the right (B->C, C->D, ...,
# J->K). Any value that was already in column K falls off the end of the
# row (the series only keeps 10 quarters, columns B..K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to place in column B for each row (the newly computed error)
$newFirstValues = @{
    2  = -0.5825945370336409
    3  = 0.09567504080935779
    4  = -0.2604190369987228
    5  = 0.8354549961584912
    6  = -0.1000793599026215
    7  = -0.3537865060796963
    8  = 0.1481773904324453
    9  = 0.157445989004155
    10 = -0.5006594565260708
    11 = 0.2803578805354692
    12 = -0.1719748578450117
    13 = 0.3058625397463315
    14 = -0.6123299526872862
    15 = 0.6883713851991116
    16 = -0.2766911554241067
}

for ($row = 2; $row -le 16; $row++) {

    # Read the existing values across the row (columns B..K => 2..11)
    # before overwriting anything.
    $oldValues = @()
    for ($col = 2; $col -le 11; $col++) {
        $oldValues += $ws.Cells.Item($row, $col).Value()
    }

    # Shift the existing values one column to the right: old column c
    # (2..10, i.e. B..J) moves to column c+1 (C..K). Whatever used to sit
    # in column K (index 10, the 10th item / oldValues[9]) is discarded.
    for ($i = 9; $i -ge 1; $i--) {
        $destCol = $i + 2
        $ws.Cells.Item($row, $destCol).Value = $oldValues[$i - 1]
    }

    # Drop the freshly computed value into column B.
    $ws.Cells.Item($row, 2).Value = $newFirstValues[$row]
}
